$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 396, pushing the existing rows 396-415 down to 397-416.
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row 396 with the new weekly price record.
$ws.Range("A396").Value = 8
$ws.Range("B396").Value = "Terminal La Palmera de La Serena"
$ws.Range("C396").Value = "Coquimbo"
$ws.Range("D396").Value = 44706
$ws.Range("E396").Value = 4
$ws.Range("F396").Value = 100114001
$ws.Range("G396").Value = "Papa"
$ws.Range("H396").Value = "Asterix"
$ws.Range("I396").Value = "1a (cosecha)"
$ws.Range("J396").Value = 2000
$ws.Range("K396").Value = 8500
$ws.Range("L396").Value = 9000
$ws.Range("M396").Value = 8750
$ws.Range("N396").Value = "`$/saco 25 kilos"
$ws.Range("O396").Value = "Región de La Araucanía"
$ws.Range("P396").Value = 350
$ws.Range("Q396").Value = 25
$ws.Range("R396").Value = "Hortaliza"
